# Update employee-count requirement figures (emp view / SF allocation
# following the previous up/down shift).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 block (anchor values feeding the C5:C10 / D5:D10 / E5:E... copy-down formulas)
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 0

# Row 6: Greeter_Up_Needed no longer follows the row above — typed directly
$ws.Range("E6").Value = 0

# Row 7: both Greeter columns become direct entries
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1

# Row 11 anchor for the C11:C20 copy-down block
$ws.Range("C11").Value = 2

# Row 21 anchors
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 2
$ws.Range("F21").Value = 1

# Row 22: Reg_Down_Needed typed directly instead of copying D21
$ws.Range("D22").Value = 2

# Update the remembered selection to match the authored file
$ws.Range("E23").Select()
